$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row with 座號 (id) 13 (sheet row 14, values 139.2/40.3/20.8/25/100/33)
# was removed from the data table. Deleting the entire worksheet row shifts
# every row below it up by one, which turns the old rows 15-18 into the new
# rows 14-17 and shrinks the used range from A1:G18 to A1:G17.
$ws.Rows(14).Delete() | Out-Null

# Reflect the resulting selection left behind after the delete (the rows
# that slid up into the deleted row's place end up selected).
$ws.Range("A14:G17").Select() | Out-Null
